# Tue, Jun 16, 2020  9:05:08 PM
# 1) The table on slide 16 (the cash-flow summary table) is switched from
#    the deck's custom "Table_0" style to PowerPoint's built-in table style
#    {EC60BE6D-03CA-4ECD-A87C-0E495E4E3825}. Table styles are not writable
#    as a plain property in this object model, so we must go through
#    Table.ApplyStyle(guid) (the setter even tells us so if we try the
#    naive `.Style = ...` assignment).
$p = $ppt.ActivePresentation

$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{EC60BE6D-03CA-4ECD-A87C-0E495E4E3825}")
    }
}

# 2) The deck's theme is swapped: the design that was using the "Integral"
#    colour palette now uses the stock "Office Theme" palette instead (the
#    12 theme colours below are exactly the Office Theme's dk1/lt1/dk2/lt2/
#    accent1-6/hlink/folHlink swatches). Every slide shares the single
#    design in this deck, so re-pointing the shared ThemeColorScheme swatch
#    by swatch reproduces the palette change everywhere at once.
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Colors(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1
$themeColors.Colors(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1
$themeColors.Colors(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2
$themeColors.Colors(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2
$themeColors.Colors(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1
$themeColors.Colors(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2
$themeColors.Colors(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3
$themeColors.Colors(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4
$themeColors.Colors(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5
$themeColors.Colors(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6
$themeColors.Colors(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink
$themeColors.Colors(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink
